$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new results add two extra tracked episode-length buckets, so the
# event-count table grows from 46 to 47 data rows. Row 47 is brand new;
# give it the same bold/bordered look as the other "length" cells in
# column A by copying the formatting down from the row above.
$ws.Range("A46").Copy($ws.Range("A47"))

# Final values (columns A:BG) for rows 42-47 after the new results were
# merged in.
$rows = @(
    @(60,0,1.15,0.02,7.09,0,0,1.14,14.86,0.03,5.99,1.12,5.98,6,14.16,14.16,14.16,0,0,0,14.16,6.97,0,0,0,0,36.99,302.58,0.01,0,0,1.2,0.03,6.94,0,0,1.2,14.82,0.03,5.96,1.05,5.96,5.96,13.94,13.94,13.94,0,0,0,13.94,6.95,0.01,0,0,0.01,37.02,303.09,0.91,0),
    @(62,0,0.42,0,0,0,0,0.42,17,0,0,0,0,0,17,17,17,0,0,0,17,0,0,0,0,0,42.99,322.17,0.55,0,0,0.3,0,0,0,0,0.3,17,0,0,0,0,0,17,17,17,0,0,0,17,0,0,0,0,0,44.21,321.19,0.61,0),
    @(66,0,1.6,0,0,0,0,1.59,13.01,0,0,0,0,0,12.96,12.96,12.96,0,0,0,12.96,0,0,0,0,0,49.31,321.52,3.24,0,0,1.59,0,0,0,0,1.56,13.03,0,0,0,0,0,12.99,12.99,12.99,0,0,0,12.99,0,0,0,0,0,47.37,323.46,3.55,0),
    @(68,0,4.76,0,0,0,0,4.76,13,0,0,0,0,0,13,13,13,0,0,0,13,0,0,0,0,0,69.40000000000001,295.02,1.32,0,0,4.78,0,0,0,0,4.77,13.01,0,0,0,0,0,13,13,13,0,0,0,13,0,0,0,0,0,68.95,295.45,1.48,0),
    @(70,0.01,0,0,16,0.01,17,0,0,0,16,0,16,16,16.93,16.93,16.93,0,0,16.93,0,16.87,0,0,0,0,1,300.18,0.5600000000000001,0,0,0,0,16,0,17,0,0,0,16,0,16,16,16.95,16.95,16.95,0,0,16.95,0,16,0,0,0,0,0.98,301.07,0.66,0),
    @(71,0.02,0,0,16.01,0.02,17,0,0,0,16,0.01,16,16,16.96,16.96,16.96,0,0,16.96,0,16.86,0,0,0,0,2.26,298.85,0.06,0,0.06,0,0,15.99,0.06,16.97,0,0,0,16,0,15.99,15.99,16.57,16.57,16.57,0,0,16.57,0,15.99,0,0,0,0,2.44,299.73,0.02,0)
)

$nRows = $rows.Count
$nCols = $rows[0].Count
$data = New-Object "object[,]" $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("A42:BG47").Value = $data

Write-Output "Updated rows 42-47 (A:BG) with the new event-count results"
